# Update LR-pair statistics in the active sheet with newly computed TPM-based values.
# Only columns E-J and M-T change for rows 2-10 (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=1; F=0.3333333333333333; G=0.047125; H=0.141375; I=0.1108387998127795; J=0.1108387998127795; M=2.761510666666667; N=8.284532; O=0.1942757069889504; P=0.1942757069889504; Q=0.1301361901666667; R=1.1712257115; S=0.02153328619543447; T=0.02153328619543447 }
    3  = @{ E=1; F=0.3333333333333333; G=0.047125; H=0.141375; I=0.1108387998127795; J=0.1108387998127795; O=0.3487673641523367; P=0.3487673641523367; Q=0.2336229100833333; R=2.10260619075; S=0.03865695605651161; T=0.03865695605651161 }
    4  = @{ E=1; F=0.3333333333333333; G=0.047125; H=0.141375; I=0.1108387998127795; J=0.1108387998127795; M=6.495363999999999; N=19.486092; O=0.4569569288587129; P=0.4569569288587129; Q=0.3060940285; R=2.7548462565; S=0.05064855756083338; T=0.05064855756083338 }
    5  = @{ I=0.5286385506557817; J=0.5286385506557816; M=2.761510666666667; N=8.284532; O=0.1942757069889504; P=0.1942757069889504; Q=0.6206762169364445; R=5.586085952428001; S=0.1027016281702661; T=0.102701628170266 }
    6  = @{ I=0.5286385506557817; J=0.5286385506557816; O=0.3487673641523367; P=0.3487673641523367; S=0.1843718739015285; T=0.1843718739015285 }
    7  = @{ I=0.5286385506557817; J=0.5286385506557816; M=6.495363999999999; N=19.486092; O=0.4569569288587129; P=0.4569569288587129; Q=1.459895847518667; R=13.139062627668; S=0.2415650485839871; T=0.2415650485839871 }
    8  = @{ E=3; F=1; G=0.1532823333333333; H=0.459847; I=0.360522649531439; J=0.360522649531439; M=2.761510666666667; N=8.284532; O=0.1942757069889504; P=0.1942757069889504; Q=0.4232907985115555; R=3.809617186604; S=0.07004079262324991; T=0.0700407926232499 }
    9  = @{ E=3; F=1; G=0.1532823333333333; H=0.459847; I=0.360522649531439; J=0.360522649531439; O=0.3487673641523367; P=0.3487673641523367; Q=0.759899517829111; R=6.839095660461999; S=0.1257385341942967; T=0.1257385341942966 }
    10 = @{ E=3; F=1; G=0.1532823333333333; H=0.459847; I=0.360522649531439; J=0.360522649531439; M=6.495363999999999; N=19.486092; O=0.4569569288587129; P=0.4569569288587129; Q=0.9956245497693332; R=8.960620947923999; S=0.1647433227138925; T=0.1647433227138924 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
